# Updated symbol list on Fri Feb 10 17:57:45 UTC 2023 with GitHub Actions
# Refresh price / volume(1h) figures for the coin rows, plus swap the
# CoinbaseStockToken / BOLO rows (48 <-> 49) to reflect the new ranking.
# Columns D (Price) and E (Volume(1h)) hold formatted numeric/percent text
# (e.g. "306.90", "-3.41%") that must stay literal strings, so each is
# forced to the "@" (Text) number format before the value is assigned -
# otherwise Excel would silently coerce them into floating point numbers
# and drop the significant trailing zeros / exact percent text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '306.90'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '-3.41%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '40.76'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '-1.96%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.029'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '-3.03%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.07602'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '-6.22%'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '4.241'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '-2.89%'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.593'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '-8.99%'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.9050'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09766'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '-12.75%'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1761'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '-4.79%'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.09189'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '-0.30%'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.04323'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-5.77%'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.1051'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-0.42%'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.001236'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-3.19%'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.005818'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '-0.38%'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.371'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '0.80%'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.419'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-6.94%'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '-2.64%'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.850'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '-7.30%'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.1349'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-2.58%'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.2724'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '6.74%'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.04159'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '-1.01%'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.001215'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '-2.40%'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.004063'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-4.61%'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0001302'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '6.40%'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0003009'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '0.59%'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02423'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '-6.00%'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05128'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '-6.39%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.007838'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '-2.82%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1303'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.007070'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '7.94%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.001950'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-6.83%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.008374'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '1.62%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.3326'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '-3.76%'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00006435'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '-4.83%'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000750'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-0.21%'
$ws.Range('B48').Value = 'BOLO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.006399'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '89.15%'
$ws.Range('B49').Value = 'CoinbaseStockToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.003003'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-27.05%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002101'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-0.21%'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0002001'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '-0.21%'
